# fix: [#64663] Fix template
#
# Adds a new leading "Generalforsamlingsdato" (general-meeting date) header
# and a new trailing "Nettoudbytte" (net dividend) header to the single
# header row of the dividend template, so the header line becomes:
#   Generalforsamlingsdato | Udbetalingsdato | Identifikation | Navn |
#   C/O | Adresse | Postnr. | Land | Bruttoudbytte | Nettoudbytte

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Generalforsamlingsdato"
$ws.Range("B1").Value = "Udbetalingsdato"
$ws.Range("C1").Value = "Identifikation"
$ws.Range("D1").Value = "Navn"
$ws.Range("E1").Value = "C/O"
$ws.Range("F1").Value = "Adresse"
$ws.Range("G1").Value = "Postnr."
$ws.Range("H1").Value = "Land"
$ws.Range("I1").Value = "Bruttoudbytte"
$ws.Range("J1").Value = "Nettoudbytte"

# Give the two new header cells the same bold header style as the rest
# of row 1 (reuses the existing header style rather than creating a new one).
$ws.Range("I1:J1").Font.Bold = $true

# The active selection moved from A2 to B2 in the updated template.
$ws.Range("B2").Select()
